# Apply the "ingestion-batch-id" StructureDefinition refresh:
#  - bump version 5.0.0 -> 6.0.0
#  - bump publication date
#  - replace the "Contact" rows with real Publisher / Jurisdiction info
#  - update the root Extension row's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value = "6.0.0"
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row gains a value, and the following duplicated "Contact" rows
# become a single "Jurisdiction" row.
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Remove the now-redundant second "Contact" row entirely (table shrinks by one row).
$ws1.Rows.Item(11).Delete()

# --- Elements sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# The root Extension row now shows the extension's own title/description
# instead of the generic placeholder text.
$ws2.Range("K2").Value = "Ingestion Batch Id"
$ws2.Range("L2").Value = "The ID generated by an ingestion service. This represents a producer-submitted message collection"
